$d = $word.ActiveDocument

function Highlight-Range([object]$rangeStart, [object]$rangeEnd, [string]$searchText) {
    $scope = $d.Range($rangeStart, $rangeEnd)
    $find = $scope.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Replacement.Highlight = $true
    $ok = $find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $true, $searchText, 2)
    return $scope.End
}

# --- Paragraph 1: highlight the word "cíclico" ---
$para1 = $d.Paragraphs(1).Range
Highlight-Range $para1.Start $para1.End "cíclico" | Out-Null

# --- Paragraph 2: highlight "cíclico", "ciclos", "ciclo" (only the first three
#     occurrences, before the unrelated later mentions of "ciclo") ---
$para2 = $d.Paragraphs(2).Range

# Find the boundary marker ("dentro de un ciclo debía enunciarse") that must stay
# un-highlighted, so the limit tracks the real text instead of a fixed offset.
$boundaryScope = $d.Range($para2.Start, $para2.End)
$boundaryScope.Find.ClearFormatting()
$boundaryScope.Find.Execute("dentro de un ciclo debía enunciarse", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$limit = $boundaryScope.Start

$pos = $para2.Start
$pos = Highlight-Range $pos $limit "cíclico"
$pos = Highlight-Range $pos $limit "ciclos"
$pos = Highlight-Range $pos $limit "ciclo"

# --- Near the end: shorten the last sentence ---
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(". Hasta entonces, serán más las preguntas que asalten a los investigadores, que las respuestas que estos puedan proporcionar…", $true, $false, $false, $false, $false, $true, 1, $false, ".", 2) | Out-Null

# --- Remove the _GoBack bookmark left over at the end of the document ---
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
} catch {
}

# --- Append a new, empty paragraph at the very end of the document ---
$endPos = $d.Content.End
$endRange = $d.Range($endPos, $endPos)
$endRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>') | Out-Null
